# Fórmulas y nombres de los datos base para los escenarios

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escenarios")

# Workbook-level defined names pointing at the base data in the "Escenarios" sheet
$wb.Names.Add('Crecimiento', '=Escenarios!$A$5')
$wb.Names.Add('PoblacionAfrica', '=Escenarios!$B$5')
$wb.Names.Add('PoblacionAmerica', '=Escenarios!$C$5')
$wb.Names.Add('PoblacionAsia', '=Escenarios!$D$5')
$wb.Names.Add('PoblacionEuropa', '=Escenarios!$E$5')
$wb.Names.Add('PoblacionOceania', '=Escenarios!$F$5')

# Forecast formulas for row 5, using the new named "Crecimiento" rate
$ws.Range("B5").Formula = '=B4*Crecimiento+B4'
$ws.Range("C5").Formula = '=C4*Crecimiento+C4'
$ws.Range("D5").Formula = '=D4*Crecimiento+D4'
$ws.Range("E5").Formula = '=E4*Crecimiento+E4'
$ws.Range("F5").Formula = '=F4*Crecimiento+F4'

# Leave the selection where the author left it
$ws.Range("F5").Select()
